$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price / 1h-volume snapshot (and, for rows 37-38, swap
# Fetch.AI <-> Dai back into rank order with their own refreshed figures).
# Cells whose new text looks like a plain number are entered with a leading
# apostrophe so Excel stores them as text (matching the sheet's existing
# inline-string cells, e.g. "1.00") instead of silently converting them to
# numeric values.

$ws.Range("D2").Value = '70.613.06'
$ws.Range("E2").Value = '  +1.33%  '
$ws.Range("D3").Value = '3.618.59'
$ws.Range("E3").Value = '  +3.04%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '''603.49'
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("D6").Value = '''198.14'
$ws.Range("E6").Value = '  +1.60%  '
$ws.Range("D7").Value = '''0.626'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").Value = '''0.217'
$ws.Range("E9").Value = '  +8.77%  '
$ws.Range("D10").Value = '''0.644'
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").Value = '''53.63'
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").Value = '''0.0000305'
$ws.Range("E12").Value = '  +1.64%  '
$ws.Range("D13").Value = '''9.52'
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").Value = '4.218.06'
$ws.Range("E14").Value = '  +3.36%  '
$ws.Range("D15").Value = '''606.34'
$ws.Range("E15").Value = '  +1.99%  '
$ws.Range("D16").Value = '''12.98'
$ws.Range("E16").Value = '  +1.70%  '
$ws.Range("D17").Value = '70.819.07'
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("D18").Value = '3.642.74'
$ws.Range("E18").Value = '  +3.37%  '
$ws.Range("D19").Value = '''19.00'
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("D21").Value = '''0.995'
$ws.Range("D22").Value = '''18.26'
$ws.Range("E22").Value = '  +1.15%  '
$ws.Range("D23").Value = '''5.31'
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("D24").Value = '''103.39'
$ws.Range("E24").Value = '  +1.52%  '
$ws.Range("D25").Value = '''4.62'
$ws.Range("E25").Value = '  -1.08%  '
$ws.Range("D26").Value = '''2.99'
$ws.Range("E26").Value = '  -5.43%  '
$ws.Range("D27").Value = '''10.56'
$ws.Range("E27").Value = '  -2.74%  '
$ws.Range("D28").Value = '''9.71'
$ws.Range("E28").Value = '  +1.90%  '
$ws.Range("D29").Value = '''33.58'
$ws.Range("E29").Value = '  +0.70%  '
$ws.Range("D30").Value = '''4.72'
$ws.Range("E30").Value = '  +13.08%  '
$ws.Range("D31").Value = '''7.17'
$ws.Range("E31").Value = '  +2.19%  '
$ws.Range("D32").Value = '''12.25'
$ws.Range("E32").Value = '  -1.08%  '
$ws.Range("D33").Value = '''0.116'
$ws.Range("E33").Value = '  +0.61%  '
$ws.Range("D34").Value = '''63.28'
$ws.Range("E34").Value = '  +0.31%  '
$ws.Range("D35").Value = '0.0₃0879'
$ws.Range("E35").Value = '  +5.02%  '
$ws.Range("D36").Value = '3.982.46'
$ws.Range("E36").Value = '  +7.35%  '
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").Value = '''1.00'
$ws.Range("E37").Value = '  +0.17%  '
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").Value = '''3.06'
$ws.Range("E38").Value = '  -0.48%  '
$ws.Range("D39").Value = '''515.41'
$ws.Range("E39").Value = '  +8.53%  '
$ws.Range("D40").Value = '''0.389'
$ws.Range("E40").Value = '  -0.48%  '
$ws.Range("D41").Value = '''36.58'
$ws.Range("E41").Value = '  +0.55%  '
$ws.Range("D42").Value = '''3.55'
$ws.Range("E42").Value = '  -2.64%  '
$ws.Range("E43").Value = '  +2.99%  '
$ws.Range("D44").Value = '''0.0461'
$ws.Range("E44").Value = '  +1.75%  '
$ws.Range("D45").Value = '''3.47'
$ws.Range("E45").Value = '  +5.77%  '
$ws.Range("D46").Value = '''2.92'
$ws.Range("E46").Value = '  +4.04%  '
$ws.Range("D47").Value = '''0.140'
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("D48").Value = '''8.58'
$ws.Range("E48").Value = '  +1.60%  '
$ws.Range("E49").Value = '  -0.16%  '
$ws.Range("D50").Value = '''0.000249'
$ws.Range("E50").Value = '  +1.79%  '
$ws.Range("D51").Value = '''1.30'
$ws.Range("E51").Value = '  +0.91%  '
